$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") '246.52'
$ws.Range("E2").Value = '1BNBBNB'
Set-TextValue $ws.Range("G2") '4'
Set-TextValue $ws.Range("D3") '22.33'
Set-TextValue $ws.Range("G3") '4'
Set-TextValue $ws.Range("D4") '5.575'
Set-TextValue $ws.Range("G4") '4'
Set-TextValue $ws.Range("D5") '0.05593'
Set-TextValue $ws.Range("G5") '4'
Set-TextValue $ws.Range("D6") '3.378'
Set-TextValue $ws.Range("G6") '4'
Set-TextValue $ws.Range("D7") '6.478'
Set-TextValue $ws.Range("G7") '4'
Set-TextValue $ws.Range("G8") '4'
Set-TextValue $ws.Range("D9") '0.8019'
Set-TextValue $ws.Range("G9") '4'
$ws.Range("B10").Value = 'One'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
Set-TextValue $ws.Range("D10") '0.0005744'
$ws.Range("E10").Value = '9OneONE'
Set-TextValue $ws.Range("G10") '4'
$ws.Range("B11").Value = 'WazirX'
$ws.Range("C11").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
Set-TextValue $ws.Range("D11") '0.1418'
$ws.Range("E11").Value = '10WazirXWRX'
Set-TextValue $ws.Range("G11") '4'
$ws.Range("B12").Value = 'MandalaExchangeToken'
$ws.Range("C12").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
Set-TextValue $ws.Range("D12") '0.07480'
$ws.Range("E12").Value = '11MandalaExchangeTokenMDX'
Set-TextValue $ws.Range("G12") '4'
$ws.Range("B13").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C13").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
Set-TextValue $ws.Range("D13") '0.03280'
$ws.Range("E13").Value = '12LiechtensteinCryptoassetsExchangeLCX'
Set-TextValue $ws.Range("G13") '4'
$ws.Range("B14").Value = 'BitrueCoin'
$ws.Range("C14").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
Set-TextValue $ws.Range("D14") '0.02989'
$ws.Range("E14").Value = '13BitrueCoinBTR'
Set-TextValue $ws.Range("G14") '4'
$ws.Range("B15").Value = 'BitMartToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
Set-TextValue $ws.Range("D15") '0.09256'
$ws.Range("E15").Value = '14BitMartTokenBMX'
Set-TextValue $ws.Range("G15") '4'
$ws.Range("B16").Value = 'BitForexToken'
$ws.Range("C16").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
Set-TextValue $ws.Range("D16") '0.001665'
$ws.Range("E16").Value = '15BitForexTokenBF'
Set-TextValue $ws.Range("G16") '4'
$ws.Range("B17").Value = 'MCDex'
$ws.Range("C17").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
Set-TextValue $ws.Range("D17") '3.256'
$ws.Range("E17").Value = '16MCDexMCB'
Set-TextValue $ws.Range("G17") '4'
$ws.Range("B18").Value = 'CoinExToken'
$ws.Range("C18").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
Set-TextValue $ws.Range("D18") '0.04713'
$ws.Range("E18").Value = '17CoinExTokenCET'
Set-TextValue $ws.Range("G18") '4'
Set-TextValue $ws.Range("D19") '0.006259'
Set-TextValue $ws.Range("G19") '4'
Set-TextValue $ws.Range("D20") '0.001048'
Set-TextValue $ws.Range("G20") '4'
Set-TextValue $ws.Range("D21") '0.003814'
$ws.Range("E21").Value = '20HotbitTokenHTBWorstin24h'
Set-TextValue $ws.Range("G21") '4'
Set-TextValue $ws.Range("D22") '0.0001499'
Set-TextValue $ws.Range("G22") '4'
Set-TextValue $ws.Range("D23") '0.0004773'
Set-TextValue $ws.Range("G23") '4'
Set-TextValue $ws.Range("D24") '3.979'
Set-TextValue $ws.Range("G24") '4'
Set-TextValue $ws.Range("D25") '2.138'
Set-TextValue $ws.Range("G25") '4'
Set-TextValue $ws.Range("G26") '4'
Set-TextValue $ws.Range("G27") '4'
Set-TextValue $ws.Range("G28") '4'
Set-TextValue $ws.Range("G29") '4'
Set-TextValue $ws.Range("G30") '4'
Set-TextValue $ws.Range("G31") '4'
Set-TextValue $ws.Range("G32") '4'
Set-TextValue $ws.Range("G33") '4'
Set-TextValue $ws.Range("G34") '4'
Set-TextValue $ws.Range("G35") '4'
Set-TextValue $ws.Range("G36") '4'
Set-TextValue $ws.Range("G37") '4'
Set-TextValue $ws.Range("G38") '4'
Set-TextValue $ws.Range("G39") '4'
Set-TextValue $ws.Range("D40") '0.04191'
Set-TextValue $ws.Range("G40") '4'
Set-TextValue $ws.Range("D41") '0.007014'
Set-TextValue $ws.Range("G41") '4'
$ws.Range("B42").Value = 'BKEXToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
Set-TextValue $ws.Range("D42") '0.1046'
$ws.Range("E42").Value = '41BKEXTokenBKK'
Set-TextValue $ws.Range("G42") '4'
$ws.Range("B43").Value = 'CEJI'
$ws.Range("C43").Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
Set-TextValue $ws.Range("D43") '0.003297'
$ws.Range("E43").Value = '42CEJICEJIBestin24h'
Set-TextValue $ws.Range("G43") '4'
Set-TextValue $ws.Range("D44") '0.009009'
Set-TextValue $ws.Range("G44") '4'
Set-TextValue $ws.Range("D45") '0.00005537'
Set-TextValue $ws.Range("G45") '4'
Set-TextValue $ws.Range("D46") '0.00000000749'
Set-TextValue $ws.Range("G46") '4'
Set-TextValue $ws.Range("D47") '0.6794'
Set-TextValue $ws.Range("G47") '4'
Set-TextValue $ws.Range("D48") '0.03032'
Set-TextValue $ws.Range("G48") '4'
Set-TextValue $ws.Range("D49") '0.00002098'
Set-TextValue $ws.Range("G49") '4'
Set-TextValue $ws.Range("D50") '0.01009'
Set-TextValue $ws.Range("G50") '4'
Set-TextValue $ws.Range("G51") '4'
